# Update hello-world local search stats: re-sort the Cloud Balancing
# benchmark results (rows 3:16) descending by score (column E), which is
# what Excel's Data > Sort does. Using the Worksheet.Sort object (rather
# than Range.Sort) also records the persisted <sortState>/<sortCondition>
# the same way the UI-driven sort does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A3:H16")
$keyRange = $ws.Range("E3:E16")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 2)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Move the active selection to E80, matching the recorded cursor position
# after the sort/review pass.
$ws.Range("E80").Select()
